# Daily attendance processing - 2026-02-07 11:01:39 UTC
# Clear the "Recorded By" column (G) contents for all data rows and
# shrink the column width now that the long name lists are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the recorded-by names from every data row (2-129). Rows that were
# already blank simply stay blank.
$ws.Range("G2:G129").ClearContents()

# Column G no longer needs to be wide enough for the long name lists.
# (ColumnWidth is in character units; the engine's stored XML "width"
# adds a fixed ~0.8333 padding on top of it, so back that off here so
# the saved column width attribute comes out to exactly 13.)
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666
